# Refresh the crypto price / 1h-volume figures (cryptos.xlsx data refresh),
# including restoring the Uniswap/Polkadot row ordering, per the data diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.784.66"
$ws.Range("E2").Value = "  +0.79%  "
$ws.Range("D3").Value = "3.391.69"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'579.80"
$ws.Range("E5").Value = "  +0.84%  "
$ws.Range("D6").Value = "'137.59"
$ws.Range("E6").Value = "  +1.38%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.389.85"
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("E9").Value = "  -0.64%  "
$ws.Range("E10").Value = "  -1.05%  "
$ws.Range("D11").Value = "'0.126"
$ws.Range("E11").Value = "  +2.90%  "
$ws.Range("E12").Value = "  +0.70%  "
$ws.Range("D13").Value = "3.970.73"
$ws.Range("E13").Value = "  +0.24%  "
$ws.Range("E14").Value = "  +1.78%  "
$ws.Range("E15").Value = "  +1.83%  "
$ws.Range("D16").Value = "3.394.77"
$ws.Range("E16").Value = "  +0.33%  "
$ws.Range("D17").Value = "'25.46"
$ws.Range("E17").Value = "  +1.03%  "
$ws.Range("D18").Value = "61.889.26"
$ws.Range("E18").Value = "  +0.83%  "
$ws.Range("D19").Value = "'14.15"
$ws.Range("E19").Value = "  +0.50%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'9.50"
$ws.Range("E20").Value = "  +0.94%  "
$ws.Range("B21").Value = "Polkadot"
$ws.Range("C21").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D21").Value = "'5.82"
$ws.Range("E21").Value = "  +0.19%  "
$ws.Range("D22").Value = "'382.92"
$ws.Range("E22").Value = "  +1.61%  "
$ws.Range("E23").Value = "  -0.83%  "
$ws.Range("D24").Value = "3.532.34"
$ws.Range("E24").Value = "  +0.52%  "
$ws.Range("E25").Value = "  +8.18%  "
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").Value = "'71.35"
$ws.Range("E27").Value = "  +0.94%  "
$ws.Range("D28").Value = "'1.73"
$ws.Range("E28").Value = "  +1.53%  "
$ws.Range("D29").Value = "'7.65"
$ws.Range("E29").Value = "  -2.13%  "
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("E31").Value = "  +3.35%  "
$ws.Range("E32").Value = "  +0.87%  "
$ws.Range("D33").Value = "'2.18"
$ws.Range("E33").Value = "  +0.83%  "
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("E35").Value = "  +0.24%  "
$ws.Range("D36").Value = "3.422.71"
$ws.Range("E36").Value = "  +0.31%  "
$ws.Range("D37").Value = "'5.40"
$ws.Range("E37").Value = "  -3.80%  "
$ws.Range("D38").Value = "'1.58"
$ws.Range("E38").Value = "  +0.94%  "
$ws.Range("E39").Value = "  -1.18%  "
$ws.Range("D40").Value = "'165.55"
$ws.Range("E40").Value = "  +2.11%  "
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("D42").Value = "'1.76"
$ws.Range("E42").Value = "  +8.63%  "
$ws.Range("D43").Value = "'0.786"
$ws.Range("E43").Value = "  +3.28%  "
$ws.Range("E44").Value = "  +1.54%  "
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("D46").Value = "'25.16"
$ws.Range("E46").Value = "  +6.10%  "
$ws.Range("D47").Value = "'4.43"
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("E48").Value = "  -0.64%  "
$ws.Range("E49").Value = "  -1.00%  "
$ws.Range("D50").Value = "'23.01"
$ws.Range("E50").Value = "  -0.21%  "
$ws.Range("D51").Value = "2.350.49"
$ws.Range("E51").Value = "  +6.86%  "
